# Add data for 2022-06-07: the "through May 29" rolling window becomes
# "through May 30", so every neighborhood's May-<year> column picks up
# the carjacking(s) that happened on May 30 of that year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / update the header title for the current ("May 2022")
# rolling column to reflect the new cutoff date.
$ws.Name = "Through 2022-05-30"
$ws.Range("B1").Value = "May 2022 (through May 30)"

# Update / insert the May-<year> counts, one per affected neighborhood row.
$ws.Range("L2").Value = 6     # Englewood       - May 2020
$ws.Range("B3").Value = 4     # Austin          - May 2022
$ws.Range("L3").Value = 6     # Austin          - May 2020
$ws.Range("AA3").Value = 5    # Austin          - May 2017
$ws.Range("B4").Value = 5     # Humboldt Park   - May 2022
$ws.Range("B5").Value = 4     # Garfield Park   - May 2022
$ws.Range("V5").Value = 6     # Garfield Park   - May 2018
$ws.Range("L14").Value = 1    # Lincoln Park    - May 2020
$ws.Range("AF21").Value = 1   # Chatham         - May 2016
$ws.Range("V23").Value = 2    # Grand Crossing  - May 2018
$ws.Range("AA23").Value = 5   # Grand Crossing  - May 2017
$ws.Range("G27").Value = 1    # Wicker Park     - May 2021
$ws.Range("Q41").Value = 1    # Morgan Park     - May 2019
$ws.Range("G42").Value = 1    # Fuller Park     - May 2021
$ws.Range("L54").Value = 1    # Bridgeport      - May 2020
$ws.Range("L56").Value = 2    # Calumet Heights - May 2020
$ws.Range("AK57").Value = 1   # Chinatown       - May 2015
$ws.Range("AA82").Value = 1   # Portage Park    - May 2017
$ws.Range("Q90").Value = 1    # Ukrainian Village - May 2019
